# Hungry Dragon - Missions content: add "kill_chain" mission type
# - rename [canBeDuringOneRun] column to [singleRunChance] in the
#   missionTypeDefinitions table (Table13)
# - convert the existing boolean flags in that column to the new
#   numeric "chance" semantics (TRUE -> 0.3, FALSE -> 0)
# - insert a new row for the "kill_chain" mission type at the end of
#   the table

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("missions")
$lo = $ws.ListObjects.Item("Table13")

# 1. Rename the column header.
$lo.HeaderRowRange.Item(6).Value = "[singleRunChance]"

# 2. Convert existing boolean values to the new numeric scale.
$ws.Range("G36").Value = 0.3
$ws.Range("G37").Value = 0.3
$ws.Range("G38").Value = 0.3
$ws.Range("G39").Value = 0
$ws.Range("G40").Value = 0
$ws.Range("G41").Value = 0
$ws.Range("G42").Value = 0.3
$ws.Range("G43").Value = 0

# 3. Insert a physical row right after the table (row 44) so
#    everything below shifts down, then grow the table into it.
$ws.Rows("44:44").Insert()
$lo.Resize($ws.Range("B35:I44"))

# 4. Populate the new "kill_chain" row, copying formats from the
#    existing rows so the new row matches the rest of the table.
$ws.Range("B43:G43").Copy()
$ws.Range("B44:G44").PasteSpecial(-4122)
$ws.Range("H39:I39").Copy()
$ws.Range("H44:I44").PasteSpecial(-4122)

$ws.Range("B44").Value = "<Definition>"
$ws.Range("C44").Value = "kill_chain"
$ws.Range("D44").Value = 0
$ws.Range("E44").Value = 7
$ws.Range("F44").Value = 1
$ws.Range("G44").Value = 1
$ws.Range("H44").Value = "TID_MISSION_OBJECTIVE_KILL_CHAIN_DESC_SINGLE_RUN"

# 5. Other tables located below row 44 need their reference ranges
#    shifted down by one row to account for the inserted row.
$loDiff = $ws.ListObjects.Item("missionDifficultyDefinitions")
$loDiff.Resize($ws.Range("B48:K51"))

$loDragon = $ws.ListObjects.Item("Table13303132")
$loDragon.Resize($ws.Range("B56:E66"))

$loDifficultyMod = $ws.ListObjects.Item("Table1330313234")
$loDifficultyMod.Resize($ws.Range("B70:D73"))

$loOtherMod = $ws.ListObjects.Item("Table133031323435")
$loOtherMod.Resize($ws.Range("B77:D78"))

# 6. The "highlight duplicates" conditional formatting on the
#    missionDifficultyDefinitions table also needs to shift down.
$cf = $ws.Range("C48:D50").FormatConditions.Item(1)
$cf.ModifyAppliesToRange($ws.Range("C49:D51"))
